$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D2:E51 is formatted as Text before writing, so numeric-looking
# strings (e.g. "0.999", "72.90") are preserved exactly as text, matching
# the inline-string cells in the source workbook.
$numRng = $ws.Range("D2:E51")
$numRng.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '49.897.74'
$ws.Range("E2").Value = '  +3.31%  '

# Row 3
$ws.Range("D3").Value = '2.602.18'
$ws.Range("E3").Value = '  +3.60%  '

# Row 4
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.10%  '

# Row 5
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '323.62'
$ws.Range("E5").Value = '  +0.70%  '

# Row 6
$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D6").Value = '109.41'
$ws.Range("E6").Value = '  +0.60%  '

# Row 7
$ws.Range("D7").Value = '0.531'
$ws.Range("E7").Value = '  +0.54%  '

# Row 8
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.06%  '

# Row 9
$ws.Range("D9").Value = '0.560'
$ws.Range("E9").Value = '  +2.87%  '

# Row 10
$ws.Range("D10").Value = '40.68'
$ws.Range("E10").Value = '  +1.82%  '

# Row 11
$ws.Range("D11").Value = '20.73'
$ws.Range("E11").Value = '  +3.10%  '

# Row 12
$ws.Range("D12").Value = '0.0821'
$ws.Range("E12").Value = '  +0.39%  '

# Row 13
$ws.Range("E13").Value = '  +0.44%  '

# Row 14
$ws.Range("E14").Value = '  +1.50%  '

# Row 15
$ws.Range("D15").Value = '3.009.02'
$ws.Range("E15").Value = '  +3.47%  '

# Row 16
$ws.Range("D16").Value = '2.575.94'
$ws.Range("E16").Value = '  +2.59%  '

# Row 17
$ws.Range("D17").Value = '0.868'
$ws.Range("E17").Value = '  +2.65%  '

# Row 18
$ws.Range("D18").Value = '49.856.14'
$ws.Range("E18").Value = '  +3.57%  '

# Row 19
$ws.Range("E19").Value = '  +11.57%  '

# Row 20
$ws.Range("D20").Value = '13.33'
$ws.Range("E20").Value = '  +1.50%  '

# Row 21
$ws.Range("D21").Value = '6.75'
$ws.Range("E21").Value = '  -0.05%  '

# Row 22
$ws.Range("D22").Value = '0.0₃0950'
$ws.Range("E22").Value = '  +0.04%  '

# Row 23
$ws.Range("D23").Value = '282.86'
$ws.Range("E23").Value = '  +1.55%  '

# Row 24
$ws.Range("D24").Value = '72.90'
$ws.Range("E24").Value = '  +0.87%  '

# Row 25
$ws.Range("D25").Value = '2.56'
$ws.Range("E25").Value = '  -0.39%  '

# Row 26
$ws.Range("D26").Value = '26.65'
$ws.Range("E26").Value = '  +3.04%  '

# Row 27
$ws.Range("D27").Value = '0.998'
$ws.Range("E27").Value = '  -0.17%  '

# Row 28
$ws.Range("E28").Value = '  +4.74%  '

# Row 29
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '2.23'
$ws.Range("E29").Value = '  -7.15%  '

# Row 30
$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D30").Value = '9.96'
$ws.Range("E30").Value = '  +1.24%  '

# Row 31
$ws.Range("D31").Value = '35.86'
$ws.Range("E31").Value = '  +1.05%  '

# Row 32
$ws.Range("D32").Value = '49.43'
$ws.Range("E32").Value = '  +0.46%  '

# Row 33
$ws.Range("D33").Value = '19.79'
$ws.Range("E33").Value = '  +1.46%  '

# Row 34
$ws.Range("D34").Value = '5.43'
$ws.Range("E34").Value = '  +0.99%  '

# Row 35
$ws.Range("D35").Value = '1.01'
$ws.Range("E35").Value = '  +0.02%  '

# Row 36
$ws.Range("D36").Value = '0.0792'
$ws.Range("E36").Value = '  +0.86%  '

# Row 37
$ws.Range("D37").Value = '2.05'
$ws.Range("E37").Value = '  +4.58%  '

# Row 38
$ws.Range("D38").Value = '4.74'
$ws.Range("E38").Value = '  +1.85%  '

# Row 39
$ws.Range("D39").Value = '3.06'
$ws.Range("E39").Value = '  +3.33%  '

# Row 40
$ws.Range("D40").Value = '124.47'
$ws.Range("E40").Value = '  +1.70%  '

# Row 41
$ws.Range("D41").Value = '22.78'
$ws.Range("E41").Value = '  +5.81%  '

# Row 42
$ws.Range("D42").Value = '0.112'
$ws.Range("E42").Value = '  +0.42%  '

# Row 43
$ws.Range("E43").Value = '  +0.37%  '

# Row 44
$ws.Range("E44").Value = '  +2.54%  '

# Row 45
$ws.Range("D45").Value = '3.36'
$ws.Range("E45").Value = '  +5.60%  '

# Row 46
$ws.Range("D46").Value = '2.041.36'
$ws.Range("E46").Value = '  +2.01%  '

# Row 47
$ws.Range("D47").Value = '2.03'
$ws.Range("E47").Value = '  +9.24%  '

# Row 48
$ws.Range("D48").Value = '2.15'
$ws.Range("E48").Value = '  +8.47%  '

# Row 49
$ws.Range("D49").Value = '9.18'
$ws.Range("E49").Value = '  +1.86%  '

# Row 50
$ws.Range("D50").Value = '5.38'
$ws.Range("E50").Value = '  +2.78%  '

# Row 51
$ws.Range("D51").Value = '81.55'
$ws.Range("E51").Value = '  +1.50%  '

# Restore default style on the numeric range (NumberFormat change above
# bumps the cell style index; reset back to Normal/General to avoid
# leaving a stray explicit style on these cells).
$numRng.Style = "Normal"
